$wb = $excel.ActiveWorkbook

# --- Sheet "Data": add new radar entry (MOTR) ---
$data = $wb.Worksheets.Item("Data")

$data.Cells.Item(27, 1).Value = "MOTR"          # Name
$data.Cells.Item(27, 2).Value = "ISRO"          # Network
$data.Cells.Item(27, 3).Value = 1.35            # Frequency (MHz)
$data.Cells.Item(27, 4).Value = "L"             # Band
$data.Cells.Item(27, 5).Value = 830             # Peak Power (KW)
$data.Cells.Item(27, 7).Value = "Phased Array"  # Antenna Shape
$data.Cells.Item(27, 9).Value = 4               # Sources

# Widen the "Peak Power (KW)" column (column E) to match column B's width
$data.Columns.Item(5).ColumnWidth = 15.67

# Update the selected cell to reflect the newly added row
$data.Range("F27").Select()

# --- Sheet "Sources": add reference for the new radar entry ---
$sources = $wb.Worksheets.Item("Sources")

$sources.Cells.Item(5, 1).Value = 4
$sources.Cells.Item(5, 2).Value = "http://www.indino.in/motr-indigenously-built-multi-object-tracking-radar-by-isro/"

$sources.Range("B5").Select()

# Restore "Data" as the active/selected sheet tab
$data.Activate()
$data.Range("F27").Select()
